# Weekly update: insert 3 new rows of "Espárragos" price data (Banquete,
# Primera, Segunda for "Sin especificar" / "Provincia de Linares", dated
# 44525) at the top of the existing block (row 32), pushing all the
# existing rows down by 3. Dimension grows from A1:R111 to A1:R114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 32 (shifts rows 32:111 down to 35:114).
$ws.Range("A32:R34").EntireRow.Insert()

# Numeric columns vs text columns for this table.
$numericCols = @(1, 4, 5, 6, 10, 11, 12, 13, 16, 17)   # A D E F J K L M P Q

function Set-DataRow($RowIndex, $Values) {
    for ($col = 1; $col -le $Values.Length; $col++) {
        $value = $Values[$col - 1]
        if ($numericCols -contains $col) {
            $ws.Cells.Item($RowIndex, $col).Value = [double]$value
        } else {
            $ws.Cells.Item($RowIndex, $col).Value = [string]$value
        }
    }
}

# New row 32: Banquete, "Sin especificar", Provincia de Linares.
Set-DataRow 32 @(
    9, "Vega Central Mapocho de Santiago", "Metropolitana", 44525, 13, 300000000,
    "Espárragos", "Sin especificar", "Banquete", 250, 1500, 1500, 1500, "`$/kilo",
    "Provincia de Linares", 1500, 1, "Hortaliza"
)

# New row 33: Primera, "Sin especificar", Provincia de Linares.
Set-DataRow 33 @(
    9, "Vega Central Mapocho de Santiago", "Metropolitana", 44525, 13, 300000000,
    "Espárragos", "Sin especificar", "Primera", 340, 1300, 1300, 1300, "`$/kilo",
    "Provincia de Linares", 1300, 1, "Hortaliza"
)

# New row 34: Segunda, "Sin especificar", Provincia de Linares.
Set-DataRow 34 @(
    9, "Vega Central Mapocho de Santiago", "Metropolitana", 44525, 13, 300000000,
    "Espárragos", "Sin especificar", "Segunda", 106, 1100, 1100, 1100, "`$/kilo",
    "Provincia de Linares", 1100, 1, "Hortaliza"
)
